# Add a title row at the top and a citation row at the bottom of the
# supplementary-materials worksheet, per "title and authors in supplementary
# materials".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row 1 with the long descriptive title, merged A1:H1.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).Insert()

# Build the combined (wrap + left-align) style on an out-of-the-way helper
# cell first, then copy/paste the format onto the target range. Doing it
# this way (one cell, two property writes) bakes both alignment attributes
# into a single new style, instead of two separate intermediate styles.
$fmtHelper = $ws.Range("Z1")
$fmtHelper.WrapText = $true
$fmtHelper.HorizontalAlignment = -4131
$fmtHelper.Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)
$fmtHelper.Clear()

$ws.Range("A1:H1").Merge()
$ws.Range("A1").Value = "Supporting Table 1: Subpopulations bioclimatic indices. Bio1 = annual mean temperature; bio2 = mean diurnal range, i.e., the mean of the monthly differences between maximum and minimum temperatures; bio7 = temperature annual range; i.e. the difference between the maximum temperature of the warmest month and the minimum temperature of the coldest month; snow = the number of days of snow cover, when the soil temperature is around 0 °C, calculated for the period in which the maximum temperature was < 0.5 °C and the minimum temperature was > -0.5 °C;  FDD = freezing degree days, i.e. the sum of daily mean temperatures for days in which the mean temperature was below 0 °C (Choler, 2018); and GDD = growing degree days, i.e. the sum of daily mean temperatures for days in which the soil mean temperature at five cm deep was above 5 °C (Körner, 2021). For easier interpretation of FDD, we transformed the values from negative to positive."
$ws.Rows.Item(1).RowHeight = 150.75

# ---------------------------------------------------------------------
# 2) Append a new row 20 with the citation, merged A20:H20.
# ---------------------------------------------------------------------
$fmtHelper2 = $ws.Range("Z2")
$fmtHelper2.WrapText = $true
$fmtHelper2.HorizontalAlignment = -4131
$fmtHelper2.Copy()
$ws.Range("A20:H20").PasteSpecial(-4122)
$fmtHelper2.Clear()

$ws.Range("A20:H20").Merge()
$ws.Range("A20").Value = "Functional intraspecific variation in the base water potential for seed germination along soil microclimatic gradients. Espinosa del Alba, C., Cruz-Tejada, D., Jiménez-Alfaro, B., and E. Fernández-Pascual. (2025). Functional Ecology."
$ws.Rows.Item(20).RowHeight = 46.5

# ---------------------------------------------------------------------
# 3) Match the saved selection shown in the diff.
# ---------------------------------------------------------------------
$ws.Range("A20:H20").Select()
